# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "NIT-9009757883" (SKY INDUSTRIAL SAS) account-statement sheet is
# refreshed with a new data pull:
#   - the worker record for NEY SMITH CERVANTES BOLAÑOS (CC 1143396862,
#     period 2006) is removed entirely - that debt was settled/retired,
#     dropping "Cant. Trabajadores" from 2 to 1 and "Cant. Periodos"
#     from 3 to 2.
#   - the remaining worker (IVETH CAROLINA MARRUGO PAUTT, CC 1047459833)
#     now has her two remaining "Periodo Mora" / "Valor Mora" pairs
#     re-ordered (2407/70763 then 2408/132680) and her "Salario Basico"
#     updated to 3548000 for both periods.
#   - the header "VALOR MORA" total is recalculated to match the new
#     sum of outstanding periods (203443).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the settled worker's row entirely (was row 18: CC, 1143396862,
# NEY SMITH CERVANTES BOLAÑOS, 2006, 35440, 886000). Excel shifts the
# signature block (rows 23-24 -> 22-23) up automatically.
$ws.Rows("18:18").Delete()

# Refresh the two remaining detail rows for IVETH CAROLINA MARRUGO PAUTT.
$ws.Range("E16").Value = "2407"
$ws.Range("F16").Value = 70763
$ws.Range("G16").Value = 3548000

$ws.Range("E17").Value = "2408"
$ws.Range("F17").Value = 132680
$ws.Range("G17").Value = 3548000

# Update the summary block: total Valor Mora, worker count, period count.
$ws.Range("E11").Value = 203443
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2
